$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (D) and Volumen (J) values between rows 4 and 5.
$ws.Range("D4").Value = 44692
$ws.Range("J4").Value = 120
$ws.Range("D5").Value = 44691
$ws.Range("J5").Value = 100
